# Master Data Tables - Test Data : master-reg_center_user_machine.xlsx
# Append 9 new rows (22-30) of registration-center/user/machine test data,
# mirroring the existing rows' shape (D/F/G text columns "eng"/"superadmin"/
# "now()", E boolean TRUE), then leave the sheet the way Excel leaves it
# after you highlight the remaining blank rows below your data (a very
# common "mark where the real data ends" gesture) and touch Page Setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id, usr_id, machine_id  — cr_by/lang_code/is_active/cr_dtimes match
# every other row in the sheet.
$newRows = @(
    @(10002, 110021, 10021),
    @(10003, 110022, 10022),
    @(10004, 110023, 10023),
    @(10005, 110024, 10024),
    @(10006, 110025, 10025),
    @(10007, 110026, 10026),
    @(10008, 110027, 10027),
    @(10009, 110028, 10028),
    @(10010, 110029, 10029)
)

$r = 22
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $r = $r + 1
}

# Touch page setup (printer-ready portrait orientation).
$ws.PageSetup.Orientation = 1

# Select from the row right after the new data down to the bottom of the
# sheet, as the last UI action of the edit.
$ws.Rows("31:1048576").Select()
